$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, formatted like the other header cells in row 1
# (bold font, thin box border, centered horizontally, aligned to top)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fill in the team record (Wins/Losses/Ties) for every data row (2..last)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 103
    $ws.Cells.Item($r, 31).Value = 59
    $ws.Cells.Item($r, 32).Value = 0
}
